$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the first empty row of the "Organizations" table with the
# Watershed Ecology Team's organization name and URL so the profile
# no longer shows a blank/placeholder row when no institution was supplied.
$ws.Range("A26").Value = "Watershed Ecology Team"
$ws.Range("B26").Value = "glfc-wet.github.io"

# Reflect the cursor position the author left the workbook in.
$ws.Application.Goto($ws.Range("H21"), $true)
